$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Recommandations sheet updates ---
$ws1.Range("D2").Value = 3565
$ws1.Range("E2").Value = 985
$ws1.Range("D3").Value = 3237.17
$ws1.Range("E3").Value = 100.18
$ws1.Range("D4").Value = 2605
$ws1.Range("D5").Value = 2568.08
$ws1.Range("E5").Value = 610.4
$ws1.Range("D6").Value = 2345
$ws1.Range("E6").Value = 595
$ws1.Range("E7").Value = 580
$ws1.Range("D8").Value = 2295
$ws1.Range("E8").Value = 575
$ws1.Range("E9").Value = 525
$ws1.Range("D10").Value = 1447.42
$ws1.Range("E10").Value = 362.67
$ws1.Range("D11").Value = 1402.54
$ws1.Range("E11").Value = 352.47
$ws1.Range("D12").Value = 1266.7
$ws1.Range("E12").Value = 319.46
$ws1.Range("D13").Value = 1040.23
$ws1.Range("E13").Value = 261.98
$ws1.Range("D14").Value = 857.51
$ws1.Range("E14").Value = 216.19
$ws1.Range("D15").Value = 751.25
$ws1.Range("E15").Value = 187.43
$ws1.Range("D16").Value = 553.02
$ws1.Range("E16").Value = 138.15
$ws1.Range("D17").Value = 515.33
$ws1.Range("E17").Value = 129
$ws1.Range("D18").Value = 485.29
$ws1.Range("E18").Value = 121.23
$ws1.Range("D19").Value = 476.93
$ws1.Range("E19").Value = 119.14
$ws1.Range("D20").Value = 428.36
$ws1.Range("E20").Value = 107.64
$ws1.Range("D21").Value = 418.38
$ws1.Range("E21").Value = 102.71
$ws1.Range("D22").Value = 368.52
$ws1.Range("E22").Value = 91.54000000000001
$ws1.Range("D23").Value = 28.93
$ws1.Range("E23").Value = 7.07
$ws1.Range("B24").Value = 2
$ws1.Range("D24").Value = 14.99
$ws1.Range("F24").Value = "🟡 Observer"
$ws1.Range("G24").Value = "➖ Neutre"
$ws1.Range("A25").Value = "BANK OF AFRICA ML (BOAM)"
$ws1.Range("B25").Value = 1
$ws1.Range("D25").Value = 7.37
$ws1.Range("E25").Value = 7.37
$ws1.Range("A26").Value = "BERNABE CI (BNBC)"
$ws1.Range("B26").Value = 2
$ws1.Range("C26").Value = 1
$ws1.Range("D26").Value = 6.94
$ws1.Range("E26").Value = 7.2
$ws1.Range("G26").Value = "👀 À surveiller"
$ws1.Range("A27").Value = "AIR LIQUIDE CI (SIVC)"
$ws1.Range("D27").Value = 4.76
$ws1.Range("E27").Value = 4.76
$ws1.Range("A28").Value = "SAPH CI (SPHC)"
$ws1.Range("C28").Value = 0
$ws1.Range("D28").Value = 4.69
$ws1.Range("E28").Value = 4.69
$ws1.Range("G28").Value = "➖ Neutre"
$ws1.Range("A29").Value = "SICABLE CI (CABC)"
$ws1.Range("D29").Value = 4.17
$ws1.Range("E29").Value = 4.17
$ws1.Range("A30").Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws1.Range("C30").Value = 1
$ws1.Range("D30").Value = 3.58
$ws1.Range("E30").Value = 7.04
$ws1.Range("G30").Value = "👀 À surveiller"
$ws1.Range("A31").Value = "SETAO CI (STAC)"
$ws1.Range("D31").Value = 3.48
$ws1.Range("E31").Value = 3.48
$ws1.Range("A32").Value = "BANK OF AFRICA BN (BOAB)"
$ws1.Range("C32").Value = 0
$ws1.Range("D32").Value = 2.86
$ws1.Range("E32").Value = 2.86
$ws1.Range("G32").Value = "➖ Neutre"
$ws1.Range("A33").Value = "SICOR CI (SICC)"
$ws1.Range("C33").Value = 0
$ws1.Range("D33").Value = 2.8
$ws1.Range("E33").Value = 2.8
$ws1.Range("G33").Value = "➖ Neutre"
$ws1.Range("A34").Value = "FILTISAC CI (FTSC)"
$ws1.Range("D34").Value = 2.09
$ws1.Range("E34").Value = -0.8
$ws1.Range("A35").Value = "UNIWAX CI (UNXC)"
$ws1.Range("B35").Value = 1
$ws1.Range("C35").Value = 0
$ws1.Range("D35").Value = 1.75
$ws1.Range("E35").Value = 1.75
$ws1.Range("A36").Value = "NEI-CEDA CI (NEIC)"
$ws1.Range("D36").Value = 0.85
$ws1.Range("E36").Value = 1.69
$ws1.Range("A37").Value = "BANK OF AFRICA NG (BOAN)"
$ws1.Range("B37").Value = 1
$ws1.Range("D37").Value = 0.05
$ws1.Range("E37").Value = 2.24
$ws1.Range("G37").Value = "👀 À surveiller"
$ws1.Range("A38").Value = "TOTAL"
$ws1.Range("C38").Value = 4
$ws1.Range("D38").Value = 0
$ws1.Range("E38").Value = 0
$ws1.Range("A39").Value = "SOLIBRA CI (SLBC)"
$ws1.Range("B39").Value = 0
$ws1.Range("C39").Value = 1
$ws1.Range("D39").Value = -0.8100000000000001
$ws1.Range("E39").Value = -0.8100000000000001
$ws1.Range("G39").Value = "➖ Neutre"
$ws1.Range("A40").Value = "SONATEL SN (SNTS)"
$ws1.Range("D40").Value = -1.4
$ws1.Range("E40").Value = -1.4
$ws1.Range("A41").Value = "BICI CI (BICC)"
$ws1.Range("D41").Value = -2.85
$ws1.Range("E41").Value = -2.85
$ws1.Range("A42").Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws1.Range("D42").Value = -3.03
$ws1.Range("E42").Value = -3.03
$ws1.Range("A44").Value = "CFAO MOTORS CI (CFAC)"
$ws1.Range("D44").Value = -3.65
$ws1.Range("E44").Value = -3.65
$ws1.Range("A45").Value = "SUCRIVOIRE (SCRC)"
$ws1.Range("D45").Value = -3.85
$ws1.Range("E45").Value = -3.85
$ws1.Range("A46").Value = "ONATEL BF (ONTBF)"
$ws1.Range("C46").Value = 2
$ws1.Range("D46").Value = -4.21
$ws1.Range("E46").Value = -2.13
$ws1.Range("A47").Value = "ORAGROUP TOGO (ORGT)"
$ws1.Range("D47").Value = -5.21
$ws1.Range("E47").Value = -1.57
$ws1.Range("A48").Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Range("C48").Value = 1
$ws1.Range("D48").Value = -6.25
$ws1.Range("E48").Value = -6.25
$ws1.Range("A49").Value = "CIE CI (CIEC)"
$ws1.Range("B49").Value = 0
$ws1.Range("C49").Value = 2
$ws1.Range("D49").Value = -6.62
$ws1.Range("E49").Value = -2.71
$ws1.Range("F49").Value = "🟡 Observer"
$ws1.Range("G49").Value = "➖ Neutre"
$ws1.Range("A50").Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws1.Range("B50").Value = 0
$ws1.Range("C50").Value = 1
$ws1.Range("D50").Value = -7.46
$ws1.Range("E50").Value = -7.46
$ws1.Range("F50").Value = "🟡 Observer"
$ws1.Range("G50").Value = "➖ Neutre"

# --- Top_YTD sheet updates ---
$ws2.Range("B2").Value = 7154304.06
$ws2.Range("B3").Value = 956088.8
$ws2.Range("B4").Value = 318415.62
$ws2.Range("B5").Value = 302645.9
$ws2.Range("A6").Value = "NEI-CEDA CI"
$ws2.Range("B6").Value = 221536.02
$ws2.Range("A7").Value = "SETAO CI"
$ws2.Range("B7").Value = 220048.64
$ws2.Range("B8").Value = 205956.41
$ws2.Range("B10").Value = 45400.57
$ws2.Range("B11").Value = 41134.07
